# Project Sample Project2 is saved. Type: SAVE.
# The only meaningful content change is cell D8 on the "Rules" sheet,
# whose value changes from 11 to 112 (kept as a plain numeric value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 112
